$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.843.11'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '1.733.86'
$ws.Range("E3").Value = '  -0.64%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''231.32'
$ws.Range("E5").Value = '  -1.87%  '
$ws.Range("D6").Value = '''0.9999'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '''0.5161'
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("E8").Value = '  +4.28%  '
$ws.Range("D9").Value = '''39.45'
$ws.Range("E9").Value = '  -2.13%  '
$ws.Range("E10").Value = '  -0.82%  '
$ws.Range("D11").Value = '1.740.77'
$ws.Range("E11").Value = '  -0.34%  '
$ws.Range("D12").Value = '''0.07033'
$ws.Range("E12").Value = '  +1.61%  '
$ws.Range("D13").Value = '''15.19'
$ws.Range("E13").Value = '  -0.50%  '
$ws.Range("D14").Value = '''0.6406'
$ws.Range("E14").Value = '  +3.26%  '
$ws.Range("D15").Value = '''4.502'
$ws.Range("E15").Value = '  +0.84%  '
$ws.Range("D16").Value = '''76.75'
$ws.Range("E16").Value = '  -0.98%  '
$ws.Range("D17").Value = '''0.9999'
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").Value = '''1.001'
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").Value = '25.836.18'
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").Value = '''11.46'
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("D21").Value = '''0.000006618'
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = '1.964.42'
$ws.Range("E22").Value = '  -0.66%  '
$ws.Range("D23").Value = '''4.134'
$ws.Range("E23").Value = '  +2.30%  '
$ws.Range("D24").Value = '''8.720'
$ws.Range("E24").Value = '  +6.04%  '
$ws.Range("D25").Value = '''5.119'
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("D26").Value = '''139.60'
$ws.Range("E26").Value = '  +2.31%  '
$ws.Range("D27").Value = '''1.511'
$ws.Range("E27").Value = '  +2.76%  '
$ws.Range("D28").Value = '''15.03'
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("D29").Value = '''1.777'
$ws.Range("E29").Value = '  +0.62%  '
$ws.Range("D30").Value = '''101.93'
$ws.Range("E30").Value = '  -0.63%  '
$ws.Range("D31").Value = '''0.08273'
$ws.Range("E31").Value = '  +0.53%  '
$ws.Range("D32").Value = '''3.681'
$ws.Range("E32").Value = '  +0.45%  '
$ws.Range("D33").Value = '''3.425'
$ws.Range("E33").Value = '  +1.37%  '
$ws.Range("D34").Value = '''0.04478'
$ws.Range("E34").Value = '  +2.36%  '
$ws.Range("D35").Value = '''2.619'
$ws.Range("E35").Value = '  -1.12%  '
$ws.Range("D36").Value = '''0.9746'
$ws.Range("E36").Value = '  -1.35%  '
$ws.Range("D37").Value = '''0.6104'
$ws.Range("E37").Value = '  +1.83%  '
$ws.Range("D38").Value = '''2.653'
$ws.Range("E38").Value = '  +0.91%  '
$ws.Range("E39").Value = '  +1.73%  '
$ws.Range("E40").Value = '  +1.02%  '
$ws.Range("D41").Value = '''0.9994'
$ws.Range("D42").Value = '''100.48'
$ws.Range("E42").Value = '  -1.53%  '
$ws.Range("D43").Value = '''0.3808'
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").Value = '''0.7236'
$ws.Range("E44").Value = '  -3.01%  '
$ws.Range("D45").Value = '''4.960'
$ws.Range("E45").Value = '  +1.58%  '
$ws.Range("D46").Value = '''0.05378'
$ws.Range("E46").Value = '  -1.86%  '
$ws.Range("E47").Value = '  +2.66%  '
$ws.Range("D48").Value = '''6.226'
$ws.Range("E48").Value = '  +5.16%  '
$ws.Range("D49").Value = '''52.90'
$ws.Range("E49").Value = '  +1.07%  '
$ws.Range("D50").Value = '''30.03'
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("D51").Value = '''7.583'
$ws.Range("E51").Value = '  +2.45%  '
